$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A350").Value = "IMX-USD"
$ws.Range("A351").Value = "TAO-USD"
$ws.Range("A352").Value = "GRT-USD"
$ws.Range("A353").Value = "MNT-USD"
